$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update cell values
$ws.Range("B2").Value = 3.8
$ws.Range("C2").Value = 11.5
$ws.Range("B3").Value = 4.5999999999999996
$ws.Range("C3").Value = 10.5

# Update column widths
$ws.Columns.Item(1).ColumnWidth = 27
$ws.Columns.Item(2).ColumnWidth = 8.43
$ws.Columns.Item(3).ColumnWidth = 27.25
